$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet1 ("openbis-metadata") edits ---

# B3: "strain1" -> "MGP90"
$ws1.Range("B3").Value = "MGP90"

# Widen column B
$ws1.Columns.Item(2).ColumnWidth = 23

# Add new row 8, copying the formatting from row 7, then fill in its values
$ws1.Range("A7:C7").Copy() | Out-Null
$ws1.Range("A8:C8").PasteSpecial(-4122) | Out-Null
$ws1.Range("A8").Value = "Header Format"
$ws1.Range("B8").Value = "TIME::VALUE_TYPE"
$ws1.Range("C8").Value = "Must be TIME::VALUE_TYPE"

# Page setup for sheet1 (Executive paper, portrait)
$ps1 = $ws1.PageSetup
$ps1.PaperSize = 10
$ps1.Orientation = 1

# --- sheet2 ("openbis-data") ---
# No data changes needed here.

# Make sheet1 the active sheet/tab with B9 selected. This also clears
# tabSelected on sheet2, since only one sheet can be the selected tab.
$ws1.Activate()
$ws1.Range("B9").Select() | Out-Null
